$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Account" column at the front of the table; the old
# PI_Account / Number Uses / Total ($) columns shift from A:C to B:D.
$ws.Columns.Item(1).Insert()

# Give the new A1 header the same (bold, bordered, centered) style the
# other header cells already carry, then set its text.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Header row
$ws.Range("A1").Value = "Account"
$ws.Range("B1").Value = "PI"
$ws.Range("C1").Value = "Number Uses"
$ws.Range("D1").Value = "Total (`$)"

# New account-coded / consolidated PI billing rows for the Feb 2024 invoice.
$data = @(
  @("CL001", "Johnson",             1.7,  68),
  @("CL002", "Hoareau/Youngquist",  4.37, 174.67),
  @("CL003", "Palatinus",           0.7,  28),
  @("CL004", "Silverton",           0,    0),
  @("CL005", "Alexander",           1,    40),
  @("CL006", "Payne",               0.2,  8),
  @("CL007", "Shah/Rieke",          7.03, 281.33),
  @("CL008", "Shah/Rieke",          3,    120)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r++
}

# Resize columns to fit the new content.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
